$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each year block, the "B" row and "C" row (3 rows after the "A" row)
# need to have their A/B/C column contents swapped.
$pairs = @(
    @(3, 4),
    @(6, 7),
    @(9, 10),
    @(12, 13),
    @(15, 16),
    @(18, 19),
    @(21, 22),
    @(24, 25),
    @(27, 28)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $a1 = $ws.Cells.Item($r1, 1).Value()
    $b1 = $ws.Cells.Item($r1, 2).Value()
    $c1 = $ws.Cells.Item($r1, 3).Value()

    $a2 = $ws.Cells.Item($r2, 1).Value()
    $b2 = $ws.Cells.Item($r2, 2).Value()
    $c2 = $ws.Cells.Item($r2, 3).Value()

    $ws.Cells.Item($r1, 1).Value = $a2
    $ws.Cells.Item($r1, 2).Value = $b2
    $ws.Cells.Item($r1, 3).Value = $c2

    $ws.Cells.Item($r2, 1).Value = $a1
    $ws.Cells.Item($r2, 2).Value = $b1
    $ws.Cells.Item($r2, 3).Value = $c1
}
